$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 489.75
$ws.Range("I4").Value = 275.81818
$ws.Range("J4").Value = 960.4
$ws.Range("K4").Value = 275.81818
$ws.Range("L4").Value = 960.4
$ws.Range("M4").Value = -161.81818
$ws.Range("N4").Value = -1188.4
$ws.Range("H45").Value = 1460
$ws.Range("I45").Value = 1460
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4380
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4188
$ws.Range("N45").ClearContents()
$ws.Range("H55").Value = 491.47058
$ws.Range("I55").Value = 519.5454999999999
$ws.Range("J55").Value = 440
$ws.Range("K55").Value = 519.5454999999999
$ws.Range("L55").Value = 440
$ws.Range("M55").Value = -305.5454999999999
$ws.Range("N55").Value = -868
$ws.Range("H112").Value = 689672.1
$ws.Range("I112").Value = 563
$ws.Range("J112").Value = 758583.0600000001
$ws.Range("K112").Value = 1689
$ws.Range("L112").Value = 2275749.18
$ws.Range("M112").Value = -581
$ws.Range("N112").Value = -2277965.18
$ws.Range("H115").Value = 478.76923
$ws.Range("I115").Value = 393.66666
$ws.Range("K115").Value = 1180.99998
$ws.Range("M115").Value = 386.0000199999999
$ws.Range("H127").Value = 756.8333
$ws.Range("I127").Value = 536.3333
$ws.Range("K127").Value = 1608.9999
$ws.Range("M127").Value = 3351.0001
$ws.Range("H129").Value = 8334.482
$ws.Range("J129").Value = 11916.65
$ws.Range("L129").Value = 35749.95
$ws.Range("N129").Value = -45749.95
$ws.Range("H138").Value = 2957.4236
$ws.Range("I138").Value = 1426.6296
$ws.Range("J138").Value = 3670.0344
$ws.Range("K138").Value = 4279.8888
$ws.Range("L138").Value = 11010.1032
$ws.Range("M138").Value = 860.1112000000003
$ws.Range("N138").Value = -21290.1032

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6662.71
$ws.Range("I32").Value = 5114.5303
$ws.Range("J32").Value = 14221.471
$ws.Range("K32").Value = 5114.5303
$ws.Range("L32").Value = 14221.471
$ws.Range("M32").Value = -4827.5303
$ws.Range("N32").Value = -14795.471

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1980.7307
$ws.Range("I99").Value = 1752.9231
$ws.Range("J99").Value = 2208.5386
$ws.Range("K99").Value = 1752.9231
$ws.Range("L99").Value = 2208.5386
$ws.Range("M99").Value = -254.9231
$ws.Range("N99").Value = -5204.5386

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 182.85715
$ws.Range("J22").Value = 133.66667
$ws.Range("L22").Value = 133.66667
$ws.Range("N22").Value = -833.6666700000001
$ws.Range("H68").Value = 24900
$ws.Range("J68").Value = 24900
$ws.Range("L68").Value = 24900
$ws.Range("N68").Value = -26398
$ws.Range("H71").Value = 24900
$ws.Range("J71").Value = 24900
$ws.Range("L71").Value = 74700
$ws.Range("N71").Value = -82188

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 586.2941
$ws.Range("I34").Value = 53.7
$ws.Range("J34").Value = 1347.1428
$ws.Range("K34").Value = 161.1
$ws.Range("L34").Value = 4041.4284
$ws.Range("M34").Value = -77.10000000000002
$ws.Range("N34").Value = -4209.428400000001
$ws.Range("H39").Value = 2782.9443
$ws.Range("J39").Value = 2974.375
$ws.Range("L39").Value = 8923.125
$ws.Range("N39").Value = -9511.125
$ws.Range("H51").Value = 3848.9375
$ws.Range("J51").Value = 4277.357
$ws.Range("L51").Value = 12832.071
$ws.Range("N51").Value = -13752.071
$ws.Range("H55").Value = 1733
$ws.Range("J55").Value = 1861.875
$ws.Range("L55").Value = 5585.625
$ws.Range("N55").Value = -5939.625
$ws.Range("H116").Value = 4215.2856
$ws.Range("I116").Value = 2131.5
$ws.Range("J116").Value = 6993.6665
$ws.Range("K116").Value = 6394.5
$ws.Range("L116").Value = 20980.9995
$ws.Range("M116").Value = -2952.5
$ws.Range("N116").Value = -27864.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1303.875
$ws.Range("I113").Value = 1071.8334
$ws.Range("K113").Value = 1071.8334
$ws.Range("M113").Value = 1098.1666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1807.8214
$ws.Range("I7").Value = 1441
$ws.Range("J7").Value = 2174.6428
$ws.Range("K7").Value = 1441
$ws.Range("L7").Value = 2174.6428
$ws.Range("M7").Value = -1329
$ws.Range("N7").Value = -2398.6428
$ws.Range("H22").Value = 2564939.5
$ws.Range("I22").Value = 8333558
$ws.Range("J22").Value = 1109
$ws.Range("K22").Value = 8333558
$ws.Range("L22").Value = 1109
$ws.Range("M22").Value = -8333263
$ws.Range("N22").Value = -1699
$ws.Range("H27").Value = 2564939.5
$ws.Range("I27").Value = 8333558
$ws.Range("J27").Value = 1109
$ws.Range("K27").Value = 8333558
$ws.Range("L27").Value = 1109
$ws.Range("M27").Value = -8333451
$ws.Range("N27").Value = -1323
$ws.Range("H55").Value = 316.12
$ws.Range("I55").Value = 417.76923
$ws.Range("J55").Value = 206
$ws.Range("K55").Value = 417.76923
$ws.Range("L55").Value = 206
$ws.Range("M55").Value = -244.76923
$ws.Range("N55").Value = -552
$ws.Range("H68").Value = 17809412
$ws.Range("I68").Value = 56389824
$ws.Range("J68").Value = 3067.6155
$ws.Range("K68").Value = 56389824
$ws.Range("L68").Value = 3067.6155
$ws.Range("M68").Value = -56389075
$ws.Range("N68").Value = -4565.6155
$ws.Range("H71").Value = 17809412
$ws.Range("I71").Value = 56389824
$ws.Range("J71").Value = 3067.6155
$ws.Range("K71").Value = 281949120
$ws.Range("L71").Value = 15338.0775
$ws.Range("M71").Value = -281945376
$ws.Range("N71").Value = -22826.0775
$ws.Range("H82").Value = 5684449.5
$ws.Range("I82").Value = 18182820
$ws.Range("J82").Value = 3372.0908
$ws.Range("K82").Value = 18182820
$ws.Range("L82").Value = 3372.0908
$ws.Range("M82").Value = -18182459
$ws.Range("N82").Value = -4094.0908
$ws.Range("H85").Value = 5684449.5
$ws.Range("I85").Value = 18182820
$ws.Range("J85").Value = 3372.0908
$ws.Range("K85").Value = 18182820
$ws.Range("L85").Value = 3372.0908
$ws.Range("M85").Value = -18181572
$ws.Range("N85").Value = -5868.0908
$ws.Range("H122").Value = 5241.6313
$ws.Range("I122").Value = 5814.769
$ws.Range("J122").Value = 3999.8333
$ws.Range("K122").Value = 17444.307
$ws.Range("L122").Value = 11999.4999
$ws.Range("M122").Value = -14994.307
$ws.Range("N122").Value = -16899.4999
$ws.Range("H126").Value = 1807.8214
$ws.Range("I126").Value = 1441
$ws.Range("J126").Value = 2174.6428
$ws.Range("K126").Value = 4323
$ws.Range("L126").Value = 6523.928400000001
$ws.Range("M126").Value = -1853
$ws.Range("N126").Value = -11463.9284
$ws.Range("H135").Value = 43868.453
$ws.Range("J135").Value = 43868.453
$ws.Range("L135").Value = 43868.453
$ws.Range("N135").Value = -54008.453
